# Update "想去人数" (column F) values across the four worksheets to reflect
# newly-scraped attendance counts, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Map of sheet name -> list of (row, newValue) pairs for column F.
$updates = @{
    "展览" = @(
        @{Row = 3;  Value = 1289},
        @{Row = 4;  Value = 1331},
        @{Row = 7;  Value = 590},
        @{Row = 10; Value = 386},
        @{Row = 13; Value = 29908},
        @{Row = 14; Value = 5611},
        @{Row = 18; Value = 69},
        @{Row = 24; Value = 679},
        @{Row = 26; Value = 318},
        @{Row = 29; Value = 117},
        @{Row = 32; Value = 233},
        @{Row = 34; Value = 599},
        @{Row = 38; Value = 255}
    )
    "演出" = @(
        @{Row = 5;  Value = 969},
        @{Row = 9;  Value = 281},
        @{Row = 16; Value = 21}
    )
    "本地生活" = @(
        @{Row = 2; Value = 303},
        @{Row = 5; Value = 318}
    )
    "全部类型" = @(
        @{Row = 2;  Value = 303},
        @{Row = 7;  Value = 318},
        @{Row = 8;  Value = 969},
        @{Row = 9;  Value = 1289},
        @{Row = 12; Value = 590},
        @{Row = 14; Value = 386},
        @{Row = 20; Value = 281},
        @{Row = 30; Value = 69},
        @{Row = 36; Value = 679},
        @{Row = 38; Value = 318},
        @{Row = 40; Value = 117},
        @{Row = 44; Value = 233},
        @{Row = 50; Value = 255}
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($entry in $updates[$sheetName]) {
        $ws.Cells.Item($entry.Row, 6).Value = $entry.Value
    }
}
